$d = $word.ActiveDocument

# 1) First occurrence (standalone paragraph after "Sincerely,"):
#    "{% if e_signature == False %}" -> "{% if add_signature == False %}"
$d.Content.Find.Execute("{% if e_", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{% if add_", 2)

# 2) Second occurrence (inside table cell): "{%if e_signature == True %} /s/ {{user.name.full(middle='full')}}"
#    becomes "{%if add_signature == True %} {{user.signature}}".
#    A "_GoBack" bookmark sits right between the runs "{%if e" and "_signature == True %}", so the
#    replacement is split into two calls that each stay on one side of the bookmark - this keeps the
#    bookmark (and everything else around it) untouched instead of having it swallowed by the replace.
$d.Content.Find.Execute("{%if e", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{%if add", 2)

$d.Content.Find.Execute(" /s/ {{user.name.full(middle='full')}}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, " {{user.signature}}", 2)
